$wb = $excel.ActiveWorkbook

# --- Sheet: "Trends Status" ---
$ws1 = $wb.Worksheets.Item("Trends Status")

$ws1.Range("B2").Value = 1
$ws1.Range("C2").Value = 18
$ws1.Range("D2").Value = 3.3
$ws1.Range("E2").Value = 20.2

$ws1.Range("B3").Value = 4
$ws1.Range("C3").Value = 24
$ws1.Range("D3").Value = 13.3
$ws1.Range("E3").Value = 27

$ws1.Range("B4").Value = 14
$ws1.Range("C4").Value = 37
$ws1.Range("D4").Value = 46.7
$ws1.Range("E4").Value = 41.6

$ws1.Range("C5").Value = 3
$ws1.Range("D5").Value = 13.3
$ws1.Range("E5").Value = 3.4

$ws1.Range("B6").Value = 7
$ws1.Range("C6").Value = 7
$ws1.Range("D6").Value = 23.3
$ws1.Range("E6").Value = 7.9

$ws1.Range("B7").Value = 49
$ws1.Range("C7").Value = 143

# --- Sheet: "Species qualification" ---
$ws4 = $wb.Worksheets.Item("Species qualification")

$ws4.Range("C3").Value = 30
$ws4.Range("C4").Value = 89

# --- Sheet: "Interannual update - High Pri" ---
$ws5 = $wb.Worksheets.Item("Interannual update - High Pri")

$ws5.Range("B2").Value = 74
$ws5.Range("C2").Value = 71.8
$ws5.Range("D2").Value = 74
$ws5.Range("E2").Value = 81.3

$ws5.Range("B3").Value = 29
$ws5.Range("C3").Value = 28.2
$ws5.Range("D3").Value = 17
$ws5.Range("E3").Value = 18.7
